$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 8-11: PUNISHMENT mute / unmute keys ---
# Key column (A) uses a dedicated Monaco 11pt font style.
# Values are entered in the same order as the original authoring session
# (all of column A first, then column C) so shared-string indices line up.

$ws.Range("A8").Value = "PUNISHMENT_MUTE_INFORM"
$ws.Range("C8").Value = "Username that was muted"
$ws.Range("A9").Value = "PUNISHMENT_MUTE_PLAYER"
$ws.Range("C9").Value = "Username that was muted"
$ws.Range("A10").Value = "PUNISHMENT_UNMUTE_PLAYER"
$ws.Range("A11").Value = "PUNISHMENT_UNMUTE_INFORM"
$ws.Range("C10").Value = "Username that was unmuted"
$ws.Range("C11").Value = "Username that was unmuted"

$ws.Range("B8").Value = "The Command Sender"
$ws.Range("B9").Value = "The Command Sender"
$ws.Range("B10").Value = "The Command Sender"
$ws.Range("B11").Value = "The Command Sender"

$ws.Range("A8").Font.Name = "Monaco"
$ws.Range("A8").Font.Size = 11
$ws.Range("A9").Font.Name = "Monaco"
$ws.Range("A9").Font.Size = 11
$ws.Range("A10").Font.Name = "Monaco"
$ws.Range("A10").Font.Size = 11
$ws.Range("A11").Font.Name = "Monaco"
$ws.Range("A11").Font.Size = 11

# Row 12 intentionally left blank (matches the source data).

# --- New row 13: PUNISHMENT-IPPARDON key ---
$ws.Range("A13").Value = "PUNISHMENT-IPPARDON"
$ws.Range("B13").Value = "The Command Sender"
$ws.Range("C13").Value = "Serialized Username array"
$ws.Range("D13").Value = "Serialized IP array"

# --- Column width tweaks (widened to fit new content) ---
$ws.Columns.Item(1).ColumnWidth = 26.65
$ws.Columns.Item(3).ColumnWidth = 24

# --- Update the active selection to land on C13 ---
$ws.Range("C13").Select()
